$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.987076
$ws.Range("N2").Value = 11.961228
$ws.Range("O2").Value = 0.2813308272685638
$ws.Range("P2").Value = 0.2813308272685638
$ws.Range("Q2").Value = 142.4935164187013
$ws.Range("R2").Value = 1282.441647768312
$ws.Range("S2").Value = 0.005485189875617521
$ws.Range("T2").Value = 0.00548518987561752

# Row 3
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("M3").Value = 10.131229
$ws.Range("N3").Value = 30.393687
$ws.Range("O3").Value = 0.7148664925918803
$ws.Range("P3").Value = 0.7148664925918804
$ws.Range("Q3").Value = 362.0784870549554
$ws.Range("R3").Value = 3258.706383494598
$ws.Range("S3").Value = 0.01393796224059
$ws.Range("T3").Value = 0.01393796224059

# Row 4
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05389233333333333
$ws.Range("N4").Value = 0.161677
$ws.Range("O4").Value = 0.00380268013955587
$ws.Range("P4").Value = 0.00380268013955587
$ws.Range("Q4").Value = 1.926050089006444
$ws.Range("R4").Value = 17.334450801058
$ws.Range("S4").Value = 0.00007414197300814044
$ws.Range("T4").Value = 0.00007414197300814044

# Row 5
$ws.Range("D5").Value = "FAPs"
$ws.Range("H5").Value = 5067.86792
$ws.Range("I5").Value = 0.9215900675332435
$ws.Range("J5").Value = 0.9215900675332435
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.987076
$ws.Range("N5").Value = 11.961228
$ws.Range("O5").Value = 0.2813308272685638
$ws.Range("P5").Value = 0.2813308272685638
$ws.Range("Q5").Value = 6735.324851667307
$ws.Range("R5").Value = 60617.92366500576
$ws.Range("S5").Value = 0.2592716961016189
$ws.Range("T5").Value = 0.2592716961016189

# Row 6
$ws.Range("D6").Value = "MuSCs"
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 10.131229
$ws.Range("N6").Value = 30.393687
$ws.Range("O6").Value = 0.7148664925918803
$ws.Range("P6").Value = 0.7148664925918804
$ws.Range("Q6").Value = 17114.57681309123
$ws.Range("R6").Value = 154031.191317821
$ws.Range("S6").Value = 0.6588138591850038
$ws.Range("T6").Value = 0.6588138591850039

# Row 7
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05389233333333333
$ws.Range("N7").Value = 0.161677
$ws.Range("O7").Value = 0.00380268013955587
$ws.Range("P7").Value = 0.00380268013955587
$ws.Range("Q7").Value = 91.03974241131554
$ws.Range("R7").Value = 819.3576817018399
$ws.Range("S7").Value = 0.003504512246620618
$ws.Range("T7").Value = 0.003504512246620618

# Row 8
$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 93.641553
$ws.Range("H8").Value = 280.924659
$ws.Range("I8").Value = 0.05108605424341119
$ws.Range("J8").Value = 0.05108605424341119
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.987076
$ws.Range("N8").Value = 11.961228
$ws.Range("O8").Value = 0.2813308272685638
$ws.Range("P8").Value = 0.2813308272685638
$ws.Range("Q8").Value = 373.355988569028
$ws.Range("R8").Value = 3360.203897121252
$ws.Range("S8").Value = 0.01437208190218559
$ws.Range("T8").Value = 0.01437208190218559

# Row 9
$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 93.641553
$ws.Range("H9").Value = 280.924659
$ws.Range("I9").Value = 0.05108605424341119
$ws.Range("J9").Value = 0.05108605424341119
$ws.Range("M9").Value = 10.131229
$ws.Range("N9").Value = 30.393687
$ws.Range("O9").Value = 0.7148664925918803
$ws.Range("P9").Value = 0.7148664925918804
$ws.Range("Q9").Value = 948.704017358637
$ws.Range("R9").Value = 8538.336156227733
$ws.Range("S9").Value = 0.0365197084173459
$ws.Range("T9").Value = 0.03651970841734591

# Row 10
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.05389233333333333
$ws.Range("N10").Value = 0.161677
$ws.Range("O10").Value = 0.00380268013955587
$ws.Range("P10").Value = 0.00380268013955587
$ws.Range("Q10").Value = 5.046561788127
$ws.Range("R10").Value = 45.419056093143
$ws.Range("S10").Value = 0.0001942639238796936
$ws.Range("T10").Value = 0.0001942639238796936

# Row 11
$ws.Range("D11").Value = "FAPs"
$ws.Range("G11").Value = 14.34625366666667
$ws.Range("H11").Value = 43.038761
$ws.Range("I11").Value = 0.007826584134129748
$ws.Range("J11").Value = 0.007826584134129748
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.987076
$ws.Range("N11").Value = 11.961228
$ws.Range("O11").Value = 0.2813308272685638
$ws.Range("P11").Value = 0.2813308272685638
$ws.Range("Q11").Value = 57.19960368427867
$ws.Range("R11").Value = 514.796433158508
$ws.Range("S11").Value = 0.002201859389141738
$ws.Range("T11").Value = 0.002201859389141738

# Row 12
$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = 14.34625366666667
$ws.Range("H12").Value = 43.038761
$ws.Range("I12").Value = 0.007826584134129748
$ws.Range("J12").Value = 0.007826584134129748
$ws.Range("M12").Value = 10.131229
$ws.Range("N12").Value = 30.393687
$ws.Range("O12").Value = 0.7148664925918803
$ws.Range("P12").Value = 0.7148664925918804
$ws.Range("Q12").Value = 145.3451811890897
$ws.Range("R12").Value = 1308.106630701807
$ws.Range("S12").Value = 0.005594962748940591
$ws.Range("T12").Value = 0.005594962748940592

# Row 13
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 14.34625366666667
$ws.Range("H13").Value = 43.038761
$ws.Range("I13").Value = 0.007826584134129748
$ws.Range("J13").Value = 0.007826584134129748
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.05389233333333333
$ws.Range("N13").Value = 0.161677
$ws.Range("O13").Value = 0.00380268013955587
$ws.Range("P13").Value = 0.00380268013955587
$ws.Range("Q13").Value = 0.7731530846885555
$ws.Range("R13").Value = 6.958377762196999
$ws.Range("S13").Value = 0.00002976199604741827
$ws.Range("T13").Value = 0.00002976199604741827
